$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the B2:B32 duty roster - each name moves up to the next slot,
# with the previous combined "崎谷航平, Jun Seomun" entry split into
# two separate rows ("Jun Seomun" then "崎谷航平").
$ws.Range("B2").Value = "川田涼介"
$ws.Range("B3").Value = "富澤天音"
$ws.Range("B4").Value = "神山修造"
$ws.Range("B5").Value = "Ethan Virtudazo"
$ws.Range("B6").Value = "豊島亮"
$ws.Range("B7").Value = "兒島大志郎"
$ws.Range("B8").Value = "高野怜央"
$ws.Range("B9").ClearContents()
$ws.Range("B10").Value = "山口玲"
$ws.Range("B11").Value = "日高泰聖"
$ws.Range("B12").Value = "志塚惇希"
$ws.Range("B13").Value = "山口洸翔"
$ws.Range("B14").Value = "白岩詩佑介"
$ws.Range("B15").Value = "石井海成"
$ws.Range("B16").Value = "Nicholas Tristan Aryasatyo"
$ws.Range("B17").Value = "小溝賢"
$ws.Range("B18").Value = "小野文哉"
$ws.Range("B19").Value = "渡部魁"
$ws.Range("B20").Value = "Jun Seomun"
$ws.Range("B21").Value = "崎谷航平"
$ws.Range("B22").Value = "三神佳誠"
$ws.Range("B23").Value = "氏家琉貴"
$ws.Range("B24").Value = "羽賀尚生"
$ws.Range("B25").Value = "島田実"
$ws.Range("B26").Value = "足立耕平"
$ws.Range("B27").Value = "遠藤隼人"
$ws.Range("B28").Value = "Ethan Virtudazo"
$ws.Range("B29").Value = "富澤天音"
$ws.Range("B30").Value = "神山修造"
$ws.Range("B31").Value = "川田涼介"
$ws.Range("B32").ClearContents()

# "Ethan Virtudazo" rows use the Roboto font style (matches the style
# already used elsewhere in the sheet for that name). Row 8 no longer
# holds that name, so restore its font back to the default Arial.
$ws.Range("B8").Font.Name = "Arial"
$ws.Range("B5").Font.Name = "Roboto"
$ws.Range("B28").Font.Name = "Roboto"

# Restore the active selection to B2.
$ws.Range("B2").Select()
